$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 "PANADOL ADVANCE 500 MG 48 TABLETS":
#  - current balance (H15) changes from "3:0" to "3:2"
#  - sale price (P15) changes from "92.0000" to "46.0000"
#  - transactions (Q15) changes from "1:0" to "0:2"

# H15 and Q15 already use a text number format, so a plain string assignment
# is stored as text.
$ws.Range("H15").Value = "3:2"
$ws.Range("Q15").Value = "0:2"

# P15 uses a numeric (0.00) number format, so assigning a numeric-looking
# string directly would be auto-converted to a number. Temporarily switch the
# cell to text format, set the text value, then restore the original numeric
# format so the style stays the same while the stored value remains text.
$p15 = $ws.Range("P15")
$p15.NumberFormat = "@"
$p15.Value = "46.0000"
$p15.NumberFormat = "0.00"

# Update the total in P23 to reflect the new sale price for row 15
# (849.70000000000005 - 92.0000 + 46.0000 = 803.70000000000005)
$ws.Range("P23").Value = 803.70000000000005

# Update the printed timestamp
$ws.Range("A24").Value = "Wednesday, 1 October, 2025 10:59 AM"
